$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Product 1
$ws.Range("A2").Value = "TC_Search_Product1"
$ws.Range("B2").Value = "Samsung Mobile"
$ws.Range("C2").Value = 50000

# Row 3 - Product 2
$ws.Range("A3").Value = "TC_Search_Product2"
$ws.Range("B3").Value = "Apple"
$ws.Range("C3").Value = 120000

# Row 4 - new header-like value (matches A1 pattern)
$ws.Range("A4").Value = 2

# Row 5 - Buy 1
$ws.Range("A5").Value = "TC_BUY_1"
$ws.Range("B5").Value = 3
$ws.Range("C5").Value = '{"quantity":1}'
$ws.Range("D5").Value = 9

# Row 6 - Buy 2
$ws.Range("A6").Value = "TC_BUY_2"
$ws.Range("B6").Value = 5
$ws.Range("C6").Value = '{"quantity":2}'
$ws.Range("D6").Value = 8
